# edit.ps1 - applies the "Added keyboard shortcuts info." change to Keyboard.docx
#
# Summary of the change:
#  - Reorders the "space bar" shortcut paragraph to after "mouse click" /
#    "mouse double-click", bolding the key/gesture name in each of those
#    three paragraphs, and moves the "_GoBack" bookmark onto "space bar".
#  - Italicizes + underlines the "The following shortcuts..." paragraph.
#  - Bolds the leading digit in each "N = toggle the ... channel on/off"
#    paragraph, and renames channel colors (green/blue/white (panneuronal)/
#    DIC/GFP (reporter)) in place of the repeated "red".
#  - Appends seven new shortcut paragraphs (arrow keys, zoom, pan, restore,
#    toggle user IDs) before the trailing blank paragraph, each with a bold
#    key name.

$d = $word.ActiveDocument

function Set-ParaText($index, $text) {
    $p = $d.Paragraphs.Item($index)
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $r.Text = $text
    return $r
}

function Bold-Range($range, $start, $len) {
    $r = $d.Range($start, $start + $len)
    $r.Bold = 1
    return $r
}

# ---------------------------------------------------------------------
# Step 0: delete the pre-existing "_GoBack" bookmark (it originally sits
# in the "6 = toggle the red channel on/off" paragraph). Bookmark names
# must be unique, and the new "_GoBack" bookmark will be (re)created
# around "space bar" later, so the stale one has to go first.
# ---------------------------------------------------------------------

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# ---------------------------------------------------------------------
# Step 1: reorder paragraphs 3-5 ("space bar" / "mouse click" /
# "mouse double-click") so the order becomes mouse click, mouse
# double-click, space bar. Re-fetch paragraph objects each time so the
# Range offsets are never stale.
# ---------------------------------------------------------------------

Set-ParaText 3 "mouse click = select a neuron" | Out-Null
Set-ParaText 4 "mouse double-click = add or remove a neuron" | Out-Null
Set-ParaText 5 "space bar (when cursor in User ID field) = save user ID" | Out-Null

# Bold "mouse click"
$p = $d.Paragraphs.Item(3)
Bold-Range $null $p.Range.Start 11 | Out-Null   # "mouse click"

# Bold "mouse double-click"
$p = $d.Paragraphs.Item(4)
Bold-Range $null $p.Range.Start 18 | Out-Null   # "mouse double-click"

# Bookmark + bold "space bar"
$p = $d.Paragraphs.Item(5)
$spaceBarRange = $d.Range($p.Range.Start, $p.Range.Start + 9)  # "space bar"
$spaceBarRange.Bold = 1
$d.Bookmarks.Add("_GoBack", $spaceBarRange) | Out-Null

# ---------------------------------------------------------------------
# Step 2: italicize + underline "The following shortcuts..." paragraph
# (paragraph mark formatting too, so new text typed at the end would
# inherit it, matching the source diff's <w:pPr><w:rPr> block).
# ---------------------------------------------------------------------

$p = $d.Paragraphs.Item(7)
$p.Range.Italic = 1
$p.Range.Underline = 1
# paragraph mark formatting
$markRange = $d.Range($p.Range.End - 1, $p.Range.End)
$markRange.Italic = 1
$markRange.Underline = 1

# ---------------------------------------------------------------------
# Step 3: bold the leading digit for the "0".."6" shortcut paragraphs,
# and rewrite the per-channel text for paragraphs 10-14 (2..6).
# ---------------------------------------------------------------------

# Paragraph 8: "0 = toggle the statistical atlas of neuron locations on/off"
$p = $d.Paragraphs.Item(8)
$d.Range($p.Range.Start, $p.Range.Start + 1).Bold = 1   # "0"

# Paragraph 9: "1 = toggle the red channel on/off"
$p = $d.Paragraphs.Item(9)
$d.Range($p.Range.Start, $p.Range.Start + 1).Bold = 1   # "1"

# Paragraph 10: "2" -> toggle the green channel on/off
Set-ParaText 10 "2 = toggle the green channel on/off" | Out-Null
$p = $d.Paragraphs.Item(10)
$d.Range($p.Range.Start, $p.Range.Start + 1).Bold = 1   # "2"

# Paragraph 11: "3" -> toggle the blue channel on/off
Set-ParaText 11 "3 = toggle the blue channel on/off" | Out-Null
$p = $d.Paragraphs.Item(11)
$d.Range($p.Range.Start, $p.Range.Start + 1).Bold = 1   # "3"

# Paragraph 12: "4" -> toggle the white (panneuronal) channel on/off
Set-ParaText 12 "4 = toggle the white (panneuronal) channel on/off" | Out-Null
$p = $d.Paragraphs.Item(12)
$d.Range($p.Range.Start, $p.Range.Start + 1).Bold = 1   # "4"

# Paragraph 13: "5" -> toggle the DIC channel on/off
Set-ParaText 13 "5 = toggle the DIC channel on/off" | Out-Null
$p = $d.Paragraphs.Item(13)
$d.Range($p.Range.Start, $p.Range.Start + 1).Bold = 1   # "5"

# Paragraph 14: "6" -> toggle the GFP (reporter) channel on/off
# (the old zero-length "_GoBack" bookmark that used to sit here was
# already removed in Step 0, before the new one was added in Step 1).
Set-ParaText 14 "6 = toggle the GFP (reporter) channel on/off" | Out-Null
$p = $d.Paragraphs.Item(14)
$d.Range($p.Range.Start, $p.Range.Start + 1).Bold = 1   # "6"

# ---------------------------------------------------------------------
# Step 4: append the new shortcut paragraphs before the trailing blank
# paragraph (currently paragraph 15).
# ---------------------------------------------------------------------

$newParas = @(
    @{ Key = "←"; Rest = " = decrease the Z slice (go left or dorsal)" },
    @{ Key = "→"; Rest = " = increase the Z slice (go right or ventral)" },
    @{ Key = "↑"; Rest = " = zoom in" },
    @{ Key = "↓"; Rest = " = zoom out" },
    @{ Key = "p"; Rest = " = toggle image panning on/off" },
    @{ Key = "r"; Rest = " = restore the image center with no zoom" },
    @{ Key = "i"; Rest = " = toggle the user IDs on/off" }
)

$insertAfter = 14
foreach ($item in $newParas) {
    $p = $d.Paragraphs.Item($insertAfter)
    $p.Range.InsertParagraphAfter()
    $insertAfter = $insertAfter + 1

    $p = $d.Paragraphs.Item($insertAfter)
    $full = $item.Key + $item.Rest
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $r.Text = $full

    $p = $d.Paragraphs.Item($insertAfter)
    $keyLen = $item.Key.Length
    $d.Range($p.Range.Start, $p.Range.Start + $keyLen).Bold = 1
}

Write-Output "done"
